$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 10175
$ws.Range("E2").Value = 744
$ws.Range("F2").Value = 744
$ws.Range("G2").Value = 1161
$ws.Range("H2").Value = 909
$ws.Range("I2").Value = 909
$ws.Range("K2").Value = 16160
$ws.Range("L2").Value = 3401
$ws.Range("M2").Value = 12759
$ws.Range("N2").Value = 12759
$ws.Range("P2").Value = 569
$ws.Range("Q2").Value = 944
$ws.Range("R2").Value = -1922
$ws.Range("S2").Value = -164
$ws.Range("T2").Value = 358
$ws.Range("U2").Value = 585
$ws.Range("V2").Value = 702
$ws.Range("W2").Value = 7.31
$ws.Range("X2").Value = 8.94
$ws.Range("Y2").Value = 7.26
$ws.Range("Z2").Value = 5.81
$ws.Range("AA2").Value = 26.66
$ws.Range("AB2").Value = 2275.98
$ws.Range("AC2").Value = 6684
$ws.Range("AD2").Value = 21.32
$ws.Range("AE2").Value = 101908
$ws.Range("AF2").Value = 1.4
$ws.Range("AG2").Value = 1460
$ws.Range("AH2").Value = 1.02
$ws.Range("AI2").Value = 19.83
$ws.Range("AJ2").Value = 13371362

# Row 3
$ws.Range("D3").Value = 11287
$ws.Range("E3").Value = 858
$ws.Range("F3").Value = 858
$ws.Range("G3").Value = 1683
$ws.Range("H3").Value = 1260
$ws.Range("I3").Value = 1260
$ws.Range("K3").Value = 18803
$ws.Range("L3").Value = 5164
$ws.Range("M3").Value = 13639
$ws.Range("N3").Value = 13639
$ws.Range("P3").Value = 569
$ws.Range("Q3").Value = 1262
$ws.Range("R3").Value = -1178
$ws.Range("S3").Value = 574
$ws.Range("T3").Value = 813
$ws.Range("U3").Value = 450
$ws.Range("V3").Value = 1516
$ws.Range("W3").Value = 7.61
$ws.Range("X3").Value = 11.17
$ws.Range("Y3").Value = 9.55
$ws.Range("Z3").Value = 7.21
$ws.Range("AA3").Value = 37.86
$ws.Range("AB3").Value = 2439.83
$ws.Range("AC3").Value = 9263
$ws.Range("AD3").Value = 24.74
$ws.Range("AE3").Value = 109239
$ws.Range("AF3").Value = 2.1
$ws.Range("AG3").Value = 1668
$ws.Range("AH3").Value = 0.73
$ws.Range("AI3").Value = 16.3
$ws.Range("AJ3").Value = 13371362

# Row 4
$ws.Range("D4").Value = 13208
$ws.Range("E4").Value = 978
$ws.Range("F4").Value = 978
$ws.Range("G4").Value = 2049
$ws.Range("H4").Value = 1612
$ws.Range("I4").Value = 1612
$ws.Range("K4").Value = 20459
$ws.Range("L4").Value = 5494
$ws.Range("M4").Value = 14966
$ws.Range("N4").Value = 14966
$ws.Range("P4").Value = 569
$ws.Range("Q4").Value = 741
$ws.Range("R4").Value = -771
$ws.Range("S4").Value = 233
$ws.Range("T4").Value = 801
$ws.Range("U4").Value = -59
$ws.Range("V4").Value = 1966
$ws.Range("W4").Value = 7.4
$ws.Range("X4").Value = 12.21
$ws.Range("Y4").Value = 11.27
$ws.Range("Z4").Value = 8.21
$ws.Range("AA4").Value = 36.71
$ws.Range("AB4").Value = 2681.45
$ws.Range("AC4").Value = 11850
$ws.Range("AD4").Value = 14.78
$ws.Range("AE4").Value = 119916
$ws.Range("AF4").Value = 1.46
$ws.Range("AG4").Value = 1668
$ws.Range("AH4").Value = 0.95
$ws.Range("AI4").Value = 12.74
$ws.Range("AJ4").Value = 13371362

# Row 5
$ws.Range("D5").Value = 14622
$ws.Range("E5").Value = 887
$ws.Range("F5").Value = 887
$ws.Range("G5").Value = 1449
$ws.Range("H5").Value = 1096
$ws.Range("I5").Value = 1090
$ws.Range("J5").Value = 6
$ws.Range("K5").Value = 20947
$ws.Range("L5").Value = 4883
$ws.Range("M5").Value = 16063
$ws.Range("N5").Value = 16000
$ws.Range("O5").Value = 63
$ws.Range("P5").Value = 595
$ws.Range("Q5").Value = 1261
$ws.Range("R5").Value = -422
$ws.Range("S5").Value = -682
$ws.Range("T5").Value = 568
$ws.Range("U5").Value = 694
$ws.Range("V5").Value = 1338
$ws.Range("W5").Value = 6.07
$ws.Range("X5").Value = 7.5
$ws.Range("Y5").Value = 7.04
$ws.Range("Z5").Value = 5.3
$ws.Range("AA5").Value = 30.4
$ws.Range("AB5").Value = 2707.01
$ws.Range("AC5").Value = 8013
$ws.Range("AD5").Value = 25.04
$ws.Range("AE5").Value = 127239
$ws.Range("AF5").Value = 1.58
$ws.Range("AG5").Value = 1745
$ws.Range("AH5").Value = 0.87
$ws.Range("AI5").Value = 19.95
$ws.Range("AJ5").Value = 13371362

# Row 6
$ws.Range("D6").Value = 15188
$ws.Range("E6").Value = 501
$ws.Range("F6").Value = 501
$ws.Range("G6").Value = 894
$ws.Range("H6").Value = 583
$ws.Range("I6").Value = 575
$ws.Range("K6").Value = 21738
$ws.Range("L6").Value = 5221
$ws.Range("M6").Value = 16517
$ws.Range("N6").Value = 16417
$ws.Range("P6").Value = 622
$ws.Range("Q6").Value = 935
$ws.Range("R6").Value = -582
$ws.Range("S6").Value = -514
$ws.Range("T6").Value = 460
$ws.Range("U6").Value = 474
$ws.Range("V6").Value = 1186
$ws.Range("W6").Value = 3.3
$ws.Range("X6").Value = 3.84
$ws.Range("Y6").Value = 3.55
$ws.Range("Z6").Value = 2.73
$ws.Range("AA6").Value = 31.61
$ws.Range("AB6").Value = 2650.39
$ws.Range("AC6").Value = 4224
$ws.Range("AD6").Value = 46.35
$ws.Range("AE6").Value = 131152
$ws.Range("AF6").Value = 1.49
$ws.Range("AG6").Value = 1826
$ws.Range("AH6").Value = 0.93
$ws.Range("AI6").Value = 39.53
$ws.Range("AJ6").Value = 13371362

# Row 7
$ws.Range("D7").Value = 15020
$ws.Range("E7").Value = 159
$ws.Range("G7").Value = 741
$ws.Range("H7").Value = 510
$ws.Range("I7").Value = 510
$ws.Range("K7").Value = 22058
$ws.Range("L7").Value = 5191
$ws.Range("M7").Value = 16867
$ws.Range("N7").Value = 16759
$ws.Range("P7").Value = 643
$ws.Range("Q7").Value = 801
$ws.Range("R7").Value = -267
$ws.Range("S7").Value = -352
$ws.Range("T7").Value = 360
$ws.Range("U7").Value = 194
$ws.Range("W7").Value = 1.06
$ws.Range("X7").Value = 3.39
$ws.Range("Y7").Value = 3.08
$ws.Range("Z7").Value = 2.33
$ws.Range("AA7").Value = 30.78
$ws.Range("AC7").Value = 3750
$ws.Range("AD7").Value = 58.14
$ws.Range("AE7").Value = 134293
$ws.Range("AF7").Value = 1.62
$ws.Range("AG7").Value = 1888
$ws.Range("AH7").Value = 0.87
$ws.Range("AI7").Value = 49.48

# Row 8
$ws.Range("D8").Value = 16167
$ws.Range("E8").Value = 678
$ws.Range("G8").Value = 1329
$ws.Range("H8").Value = 944
$ws.Range("I8").Value = 935
$ws.Range("K8").Value = 23123
$ws.Range("L8").Value = 5490
$ws.Range("M8").Value = 17633
$ws.Range("N8").Value = 17514
$ws.Range("P8").Value = 648
$ws.Range("Q8").Value = 894
$ws.Range("R8").Value = -478
$ws.Range("S8").Value = -230
$ws.Range("T8").Value = 409
$ws.Range("U8").Value = 314
$ws.Range("W8").Value = 4.19
$ws.Range("X8").Value = 5.84
$ws.Range("Y8").Value = 5.46
$ws.Range("Z8").Value = 4.18
$ws.Range("AA8").Value = 31.14
$ws.Range("AC8").Value = 6874
$ws.Range("AD8").Value = 31.71
$ws.Range("AE8").Value = 140350
$ws.Range("AF8").Value = 1.55
$ws.Range("AG8").Value = 1734
$ws.Range("AH8").Value = 0.8
$ws.Range("AI8").Value = 24.79

# Row 9
$ws.Range("D9").Value = 17150
$ws.Range("E9").Value = 806
$ws.Range("G9").Value = 1404
$ws.Range("H9").Value = 1001
$ws.Range("I9").Value = 984
$ws.Range("K9").Value = 24320
$ws.Range("L9").Value = 5843
$ws.Range("M9").Value = 18477
$ws.Range("N9").Value = 18341
$ws.Range("P9").Value = 650
$ws.Range("Q9").Value = 1038
$ws.Range("R9").Value = -504
$ws.Range("S9").Value = -194
$ws.Range("T9").Value = 344
$ws.Range("U9").Value = 820
$ws.Range("W9").Value = 4.7
$ws.Range("X9").Value = 5.84
$ws.Range("Y9").Value = 5.49
$ws.Range("Z9").Value = 4.22
$ws.Range("AA9").Value = 31.62
$ws.Range("AC9").Value = 7230
$ws.Range("AD9").Value = 30.15
$ws.Range("AE9").Value = 146974
$ws.Range("AF9").Value = 1.48
$ws.Range("AG9").Value = 1756
$ws.Range("AH9").Value = 0.81
$ws.Range("AI9").Value = 23.86

# Clear cells removed in the diff (J and O columns for rows 2-4)
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()
